$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 21: B21/C21 revert to the "normal" styles used by the rest of the
# table (they previously carried the special "last row" styles).
$ws.Range("B21").Style = $ws.Range("B20").Style
$ws.Range("C21").Style = $ws.Range("C20").Style

# New row 22 picks up the special "last row" styles that row 21 used to have.
$ws.Range("A22").Value = 45710
$ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat
$ws.Range("A22").Style = $ws.Range("A21").Style

$ws.Range("B22").Value = 4
$ws.Range("B22").Style = "Cell Check"

$ws.Range("C22").Value = "Debugged the code for automated scanning tool and started generating reports."

$ws.Range("C21").Select()
